$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace the Chinese header/value strings with their English equivalents.
$ws.Range("A1").Value = "StringA"
$ws.Range("B1").Value = "StringB"
$ws.Range("C1").Value = "StringC"

for ($i = 0; $i -le 9; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = "StringA$i"
}
for ($i = 0; $i -le 9; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = "StringB$i"
}
for ($i = 0; $i -le 9; $i++) {
    $ws.Cells.Item($i + 2, 3).Value = "StringC$i"
}

# Update the selected range to C2:C11 with C2 as the active cell.
$ws.Range("C2:C11").Select()
